# Appends the new match row (row 83) to the "Premijer Liga BiH 2023-2024"
# results sheet, mirroring the existing rows' layout/formatting.
#   Zrinjski 3 - 1 Siroki Brijeg

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 82
$newRow = $lastRow + 1

# Copy the formatting (styles) of the last existing row onto the new row
# before filling in values, so fonts/borders/alignment/number-formats match.
$ws.Range("A" + $lastRow + ":V" + $lastRow).Copy()
$ws.Range("A" + $newRow + ":V" + $newRow).PasteSpecial(-4122)

$ws.Range("A" + $newRow).Value = 82
$ws.Range("B" + $newRow).Value = "bosnia-and-herzegovina"
$ws.Range("C" + $newRow).Value = "premijer-liga-bih"
$ws.Range("D" + $newRow).Value = "2023-2024"
$ws.Range("E" + $newRow).Value = 45252.75
$ws.Range("F" + $newRow).Value = "Zrinjski"
$ws.Range("G" + $newRow).Value = 3
$ws.Range("H" + $newRow).Value = "Siroki Brijeg"
$ws.Range("I" + $newRow).Value = 1
$ws.Range("J" + $newRow).Value = 1.31
$ws.Range("K" + $newRow).Value = "03/09/2023 08:12"
$ws.Range("L" + $newRow).Value = 1.16
$ws.Range("M" + $newRow).Value = "22/11/2023 17:59"
$ws.Range("N" + $newRow).Value = 4.72
$ws.Range("O" + $newRow).Value = "03/09/2023 08:12"
$ws.Range("P" + $newRow).Value = 6.82
$ws.Range("Q" + $newRow).Value = "22/11/2023 17:59"
$ws.Range("R" + $newRow).Value = 7.8
$ws.Range("S" + $newRow).Value = "03/09/2023 08:12"
$ws.Range("T" + $newRow).Value = 15.52
$ws.Range("U" + $newRow).Value = "22/11/2023 17:59"
$ws.Range("V" + $newRow).Value = "https://www.betexplorer.com/football/bosnia-and-herzegovina/premijer-liga-bih/zrinjski-siroki-brijeg/l6bIl11G/"
